$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1:D19").Value = "Bộ môn công nghệ phần mềm"

$ws.Range("A21").Select()
